# Appends a "Note:" callout paragraph (and a preceding spacer paragraph)
# right after the closing line "True motivation doesn't come before
# action. It comes from it.", while preserving the document's existing
# trailing empty paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the body content (the one
# right before the document's trailing empty paragraph).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$target = $lastPara.Range
$target.Collapse(1)  # wdCollapseStart -- start of the trailing empty paragraph

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$fragment = @"
<w:p $wns>
  <w:pPr>
    <w:rPr>
      <w:color w:val="C00000"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p $wns>
  <w:pPr>
    <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1" w:line="240" w:lineRule="auto"/>
    <w:jc w:val="center"/>
    <w:outlineLvl w:val="2"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="C00000"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="C00000"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>Note:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="C00000"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:color w:val="C00000"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
    </w:rPr>
    <w:t>This is a default pdf for testing my code</w:t>
  </w:r>
</w:p>
<w:p $wns/>
"@

$target.InsertXML($fragment)
